$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.930.78"
$ws.Range("E2").Value = "  +1.97%  "

# Row 3
$ws.Range("D3").Value = "2.221.43"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "292.75"
$ws.Range("E5").Value = "  -0.30%  "

# Row 6
$ws.Range("D6").Value = "86.85"
$ws.Range("E6").Value = "  +7.26%  "

# Row 7
$ws.Range("D7").Value = "0.515"
$ws.Range("E7").Value = "  +1.35%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +2.28%  "

# Row 10
$ws.Range("D10").Value = "30.80"
$ws.Range("E10").Value = "  +6.50%  "

# Row 11
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  +2.71%  "

# Row 12
$ws.Range("D12").Value = "47.60"
$ws.Range("E12").Value = "  +1.17%  "

# Row 13
$ws.Range("E13").Value = "  +1.73%  "

# Row 14
$ws.Range("D14").Value = "6.36"
$ws.Range("E14").Value = "  +2.60%  "

# Row 15
$ws.Range("D15").Value = "2.559.18"
$ws.Range("E15").Value = "  +1.34%  "

# Row 16
$ws.Range("D16").Value = "14.08"
$ws.Range("E16").Value = "  +1.46%  "

# Row 17
$ws.Range("D17").Value = "2.217.76"

# Row 18
$ws.Range("D18").Value = "0.733"
$ws.Range("E18").Value = "  +3.89%  "

# Row 19
$ws.Range("D19").Value = "39.840.76"
$ws.Range("E19").Value = "  +2.02%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0882"
$ws.Range("E20").Value = "  +2.08%  "

# Row 21
$ws.Range("D21").Value = "11.31"
$ws.Range("E21").Value = "  +11.06%  "

# Row 22
$ws.Range("D22").Value = "5.82"
$ws.Range("E22").Value = "  +2.35%  "

# Row 23
$ws.Range("D23").Value = "65.79"
$ws.Range("E23").Value = "  +1.93%  "

# Row 24
$ws.Range("D24").Value = "236.36"
$ws.Range("E24").Value = "  +5.38%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +3.76%  "

# Row 27
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").Value = "  +3.45%  "

# Row 28
$ws.Range("D28").Value = "22.83"
$ws.Range("E28").Value = "  +2.10%  "

# Row 29
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +1.43%  "

# Row 30
$ws.Range("D30").Value = "9.27"
$ws.Range("E30").Value = "  +2.83%  "

# Row 31
$ws.Range("D31").Value = "32.96"
$ws.Range("E31").Value = "  +5.24%  "

# Row 32
$ws.Range("D32").Value = "152.00"
$ws.Range("E32").Value = "  +2.18%  "

# Row 33
$ws.Range("E33").Value = "  -0.20%  "

# Row 34
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  +3.52%  "

# Row 35
$ws.Range("D35").Value = "0.0721"
$ws.Range("E35").Value = "  +4.44%  "

# Row 36
$ws.Range("E36").Value = "  +1.94%  "

# Row 37
$ws.Range("D37").Value = "2.81"
$ws.Range("E37").Value = "  +7.59%  "

# Row 38
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  +2.49%  "

# Row 39
$ws.Range("D39").Value = "15.95"
$ws.Range("E39").Value = "  +4.48%  "

# Row 40
$ws.Range("D40").Value = "0.0993"
$ws.Range("E40").Value = "  +3.88%  "

# Row 41
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +4.87%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "3.79"
$ws.Range("E42").Value = "  +6.48%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.068.46"
$ws.Range("E43").Value = "  +9.66%  "

# Row 44
$ws.Range("D44").Value = "0.0269"
$ws.Range("E44").Value = "  +4.34%  "

# Row 45
$ws.Range("D45").Value = "9.97"
$ws.Range("E45").Value = "  +12.04%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "17.83"
$ws.Range("E46").Value = "  +11.70%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("D48").Value = "2.63"
$ws.Range("E48").Value = "  +1.84%  "

# Row 49
$ws.Range("D49").Value = "2.429.54"
$ws.Range("E49").Value = "  +1.38%  "

# Row 50
$ws.Range("D50").Value = "71.41"
$ws.Range("E50").Value = "  +0.26%  "

# Row 51
$ws.Range("D51").Value = "89.40"
$ws.Range("E51").Value = "  +3.29%  "
